# Updates cryptos list figures (Price + Volume(1h) columns) for rows 2-51
# on the active worksheet, matching the scraped GitHub Actions data refresh.
#
# Column D (Price) and column E (Volume(1h)) are stored as plain text in the
# workbook (values like "67.330.28" use '.' as a thousands separator, and
# values like "597.10" / "1.00" rely on an exact trailing zero for display).
# A leading apostrophe forces Excel to keep number-looking strings as literal
# text instead of silently coercing them to numeric values (which would drop
# meaningful trailing zeros, e.g. "597.10" -> 597.1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.330.28'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '3.521.92'
$ws.Range("E3").Value = '  +0.52%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''597.10'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").Value = '''173.68'
$ws.Range("E6").Value = '  +2.76%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  +3.22%  '
$ws.Range("E9").Value = '  +8.17%  '
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").Value = '''0.437'
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").Value = '4.132.45'
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("E14").Value = '  +2.14%  '
$ws.Range("E15").Value = '  +1.80%  '
$ws.Range("D16").Value = '67.288.40'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").Value = '3.513.32'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = '''6.33'
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").Value = '''14.28'
$ws.Range("E19").Value = '  +1.57%  '
$ws.Range("D20").Value = '''397.94'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").Value = '''7.99'
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").Value = '''73.46'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '''0.539'
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("E25").Value = '  -3.56%  '
$ws.Range("D26").Value = '''10.26'
$ws.Range("E26").Value = '  +2.26%  '
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("E29").Value = '  -1.50%  '
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("D32").Value = '''24.14'
$ws.Range("E32").Value = '  +2.47%  '
$ws.Range("D33").Value = '''7.40'
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("D34").Value = '''1.64'
$ws.Range("E34").Value = '  +2.38%  '
$ws.Range("D35").Value = '''163.40'
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("D36").Value = '''0.896'
$ws.Range("E36").Value = '  -0.45%  '
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("E38").Value = '  +3.83%  '
$ws.Range("D39").Value = '''4.72'
$ws.Range("E39").Value = '  +1.16%  '
$ws.Range("D40").Value = '''27.53'
$ws.Range("E40").Value = '  +3.50%  '
$ws.Range("E41").Value = '  -0.91%  '
$ws.Range("D42").Value = '''26.40'
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").Value = '''2.62'
$ws.Range("E43").Value = '  +2.83%  '
$ws.Range("D44").Value = '2.804.37'
$ws.Range("E44").Value = '  -1.14%  '
$ws.Range("D45").Value = '''42.91'
$ws.Range("E45").Value = '  -1.36%  '
$ws.Range("E46").Value = '  -2.67%  '
$ws.Range("D47").Value = '''340.40'
$ws.Range("E47").Value = '  -2.37%  '
$ws.Range("E48").Value = '  +1.59%  '
$ws.Range("D49").Value = '''33.72'
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("D50").Value = '''6.54'
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("E51").Value = '  -0.66%  '
